$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.8745545148849487
$ws.Range("B1").Value = 1.796998381614685
$ws.Range("D1").Value = 1.901783227920532
$ws.Range("E1").Value = 1.125640988349915
